$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 0.0073333333333333332, 32.759981167585842, 16.255830592676201, 0.49621001030246104, 29.484157969007036, 14.473890698248377, 0.49090398692962323, 299.5, 267.53388567723931),
    @(3, 0.0076666666666666671, 32.773716321270662, 16.267584775096754, 0.49636070000822075, 29.456641558110285, 14.446280068969243, 0.49042522517275072, 299.5, 267.06634311906811),
    @(4, 0.0080000000000000002, 32.776510725160918, 16.268560151395505, 0.49634814052695764, 29.468268198019736, 14.543245366568563, 0.49352222766677095, 299.42885403279951, 267.07683827884216),
    @(5, 0.0083333333333333332, 32.770422117765463, 16.261398766668357, 0.49622182797128961, 29.448032252660553, 14.521508767853348, 0.4931232295339994, 299.46505562653471, 267.02984737512418),
    @(6, 0.0090000000000000011, 32.78687393164747, 16.27451097026638, 0.49637275588379404, 29.40862926909962, 14.534799489855429, 0.49423587059623703, 299.44182270258415, 266.97448582686991),
    @(7, 0.0093333333333333341, 32.76191764431006, 16.248612515834804, 0.4959603614245941, 29.357828359400678, 14.513303623197459, 0.49435889622095125, 298.90000000000003, 266.07095207366655),
    @(8, 0.0096666666666666672, 32.78492082339698, 16.268758823157672, 0.49622687548317829, 29.346994149270166, 14.513565164472201, 0.49455031376128655, 299.00567760877067, 266.06085905064111),
    @(9, 0.01, 32.807420928529716, 16.291181932318644, 0.49657002809847944, 29.344262227256216, 14.583801195120969, 0.49698987427855335, 299.06489935018368, 265.58113131378332)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
    $ws.Cells.Item($r, 8).Value = $row[8]
    $ws.Cells.Item($r, 9).Value = $row[9]
}
